# Update computed statistics (MAX/MIN/MEAN/STD DEV/STD ERROR of Length/Total Width)
# for angular.js, atom, axios, and create-react-app blocks after filtering out
# JS files with length > 1000 from the underlying dataset ("Removed all js > 1000").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 30.592414921835317
$ws.Range("F3").Value = 165.06781659292028
$ws.Range("G3").Value = 1.2400602963634322
$ws.Range("L3").Value = 730
$ws.Range("N3").Value = 18.134529147982065
$ws.Range("O3").Value = 64.150186873456377
$ws.Range("P3").Value = 2.4801895136795258
$ws.Range("C4").Value = 544
$ws.Range("E4").Value = 85.617529205937132
$ws.Range("F4").Value = 48.583826805518378
$ws.Range("G4").Value = 0.36498256238221033
$ws.Range("L4").Value = 479
$ws.Range("N4").Value = 82.139013452914796
$ws.Range("O4").Value = 54.106356974901608
$ws.Range("P4").Value = 2.0918726153871567
$ws.Range("E5").Value = 6.9038320446977819
$ws.Range("F5").Value = 4.1773709170154349
$ws.Range("G5").Value = 0.031382203534861076
$ws.Range("N5").Value = 5.6591928251121075
$ws.Range("O5").Value = 2.4793002424606447
$ws.Range("P5").Value = 0.095855285265869478
$ws.Range("L10").Value = 462
$ws.Range("N10").Value = 27.366403607666292
$ws.Range("O10").Value = 54.91756576103699
$ws.Range("P10").Value = 1.8439513932942531
$ws.Range("L11").Value = 313
$ws.Range("N11").Value = 66.096956031567075
$ws.Range("O11").Value = 32.155654089834329
$ws.Range("P11").Value = 1.0796811974376637
$ws.Range("N12").Value = 4.5839909808342725
$ws.Range("O12").Value = 3.5539794092547603
$ws.Range("P12").Value = 0.11933094980848359
$ws.Range("E17").Value = 17.486425339366516
$ws.Range("F17").Value = 65.934415147148016
$ws.Range("G17").Value = 2.2176143628772422
$ws.Range("N17").Value = 29.977682596934176
$ws.Range("O17").Value = 197.23938232098257
$ws.Range("P17").Value = 2.9614054075571348
$ws.Range("C18").Value = 124
$ws.Range("E18").Value = 66.237556561085967
$ws.Range("F18").Value = 20.9204513894251
$ws.Range("G18").Value = 0.70363092438943098
$ws.Range("L18").Value = 804
$ws.Range("N18").Value = 76.958521190261493
$ws.Range("O18").Value = 44.74458213320132
$ws.Range("P18").Value = 0.67180725232909289
$ws.Range("E19").Value = 4.9276018099547514
$ws.Range("F19").Value = 3.4906738247717946
$ws.Range("G19").Value = 0.11740406573195188
$ws.Range("N19").Value = 2.088142470694319
$ws.Range("O19").Value = 5.8113670889284741
$ws.Range("P19").Value = 0.087253436509172985
$ws.Range("E20").Value = 0.15158371040723981
$ws.Range("F20").Value = 0.35861691140604879
$ws.Range("G20").Value = 0.012061591988491745
$ws.Range("N20").Value = 2.1368349864743013
$ws.Range("O20").Value = 1.7466921696279536
$ws.Range("P20").Value = 0.026225308432856777
$ws.Range("E24").Value = 46.356000000000002
$ws.Range("F24").Value = 69.482985428088796
$ws.Range("G24").Value = 4.3944898516210049
$ws.Range("L24").Value = 16517
$ws.Range("N24").Value = 20.685010958464744
$ws.Range("O24").Value = 171.37448111493558
$ws.Range("P24").Value = 0.62648991469281545
$ws.Range("C25").Value = 232
$ws.Range("E25").Value = 78.5
$ws.Range("F25").Value = 42.85708809520311
$ws.Range("G25").Value = 2.7105202452665802
$ws.Range("L25").Value = 976
$ws.Range("N25").Value = 63.415847009140961
$ws.Range("O25").Value = 46.607155118333665
$ws.Range("P25").Value = 0.17038075006381415
$ws.Range("E26").Value = 4.5759999999999996
$ws.Range("F26").Value = 3.3015487274913875
$ws.Range("G26").Value = 0.20880827569806712
$ws.Range("N26").Value = 4.2267466723686322
$ws.Range("O26").Value = 4.5735202024716397
$ws.Range("P26").Value = 0.016719317035134778
$ws.Range("N27").Value = 0.032768482386272516
$ws.Range("O27").Value = 0.24516633822593475
$ws.Range("P27").Value = 0.00089624918086669446

# Update the active cell selection to match the latest edit location.
$ws.Range("Q21").Select()
